# daily auto push: 2026-02-14 13:52 UTC
# Two new readings for 2026/02/14 (土) were appended to the source data,
# pushing every subsequent row down by two positions (819/820 -> insert,
# old 819..860 -> new 821..862).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the old row 819, shifting existing data down.
$ws.Range("A819:D820").Insert()

# New row 819: 2026/02/14, 土, 18:00, ranking 201
$ws.Range("A819").Value = "'2026/02/14"
$ws.Range("A819").Style = "Normal"
$ws.Range("B819").Value = "土"
$ws.Range("C819").Value = 18
$ws.Range("D819").Value = 201

# New row 820: 2026/02/14, 土, 21:00, ranking 201
$ws.Range("A820").Value = "'2026/02/14"
$ws.Range("A820").Style = "Normal"
$ws.Range("B820").Value = "土"
$ws.Range("C820").Value = 21
$ws.Range("D820").Value = 201
